$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A2").Value = "'4,404"
$ws1.Range("A4").Value = "'1874,3"
$ws1.Range("A8").Value = "'26408,37"
$ws1.Range("A34").Value = "'6651,0"

$ws2 = $wb.Worksheets.Item("data")
$ws2.Range("A1").Value = "-"
$ws2.Range("A2").Value = "'"
$ws2.Range("A3").Value = "'"
